# Update the student export sheet:
#  - normalise the CNE identifiers to the new 8-digit numbering scheme
#  - re-case the header row (FIRSTNAME/LASTNAME/... -> FirstName/LastName/...,
#    with the "ClasseName" header, matching the converted XML export headers)
#  - move the active selection, matching where the author left off editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CNE numbering scheme: row 2 holds the literal seed value, rows 3:11
# keep their existing "=previous+1" formulas, so updating A2 ripples through
# the whole column automatically.
$ws.Range("A2").Value = 21000001

# Re-cased header row (B1:G1); A1 ("CNE") is unchanged.
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "DateofBirth"
$ws.Range("E1").Value = "ClasseName"
$ws.Range("F1").Value = "Phone"
$ws.Range("G1").Value = "Email"

# Leave the selection where the author left it.
$ws.Range("C14").Select() | Out-Null
